# Updates coin "Price" (D) and "Volume(1h)" (E) columns with freshly scraped values.
# Cells whose new text would otherwise be auto-parsed by Excel as a number (plain
# decimals with a single "." and no thousands separators) are first formatted as
# Text ("@") so the value round-trips as the exact literal string from the source feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '59.550.57'
$ws.Range("E2").Value = '  +1.00%  '

# Row 3
$ws.Range("D3").Value = '2.603.50'
$ws.Range("E3").Value = '  +0.72%  '

# Row 4
$ws.Range("E4").Value = '  +0.45%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '539.86'
$ws.Range("E5").Value = '  +3.45%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.69'
$ws.Range("E6").Value = '  +2.01%  '

# Row 7
$ws.Range("E7").Value = '  +0.13%  '

# Row 8
$ws.Range("E8").Value = '  +0.44%  '

# Row 10
$ws.Range("E10").Value = '  +1.45%  '

# Row 11
$ws.Range("E11").Value = '  +1.55%  '

# Row 12
$ws.Range("E12").Value = '  -0.72%  '

# Row 13
$ws.Range("D13").Value = '3.061.91'
$ws.Range("E13").Value = '  +0.53%  '

# Row 14
$ws.Range("D14").Value = '59.464.70'
$ws.Range("E14").Value = '  +0.97%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.75'
$ws.Range("E15").Value = '  +1.07%  '

# Row 16
$ws.Range("D16").Value = '2.603.75'
$ws.Range("E16").Value = '  +1.42%  '

# Row 17
$ws.Range("E17").Value = '  +0.58%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '341.10'
$ws.Range("E18").Value = '  +0.29%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.37'
$ws.Range("E19").Value = '  +1.70%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.08'
$ws.Range("E20").Value = '  +0.28%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.33'
$ws.Range("E21").Value = '  -1.69%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.05%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.26'
$ws.Range("E23").Value = '  +1.71%  '

# Row 24
$ws.Range("E24").Value = '  +1.48%  '

# Row 25
$ws.Range("E25").Value = '  -1.42%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.994'
$ws.Range("E26").Value = '  -0.49%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.22'
$ws.Range("E27").Value = '  +2.99%  '

# Row 28
$ws.Range("D28").Value = '0.0₃0744'
$ws.Range("E28").Value = '  +3.25%  '

# Row 29
$ws.Range("E29").Value = '  +0.04%  '

# Row 30
$ws.Range("E30").Value = '  +6.34%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.84'
$ws.Range("E31").Value = '  -0.66%  '

# Row 32
$ws.Range("E32").Value = '  +0.80%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '149.63'
$ws.Range("E33").Value = '  +0.11%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.00'
$ws.Range("E34").Value = '  +0.91%  '

# Row 35
$ws.Range("E35").Value = '  +0.46%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.846'
$ws.Range("E36").Value = '  +4.31%  '

# Row 37
$ws.Range("E37").Value = '  -0.71%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.824'
$ws.Range("E38").Value = '  +0.23%  '

# Row 39
$ws.Range("E39").Value = '  +0.32%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.16%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '273.11'
$ws.Range("E41").Value = '  +0.28%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.600'
$ws.Range("E42").Value = '  +1.37%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.71'
$ws.Range("E43").Value = '  -0.50%  '

# Row 44
$ws.Range("E44").Value = '  +0.10%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0522'
$ws.Range("E45").Value = '  +1.44%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '18.55'
$ws.Range("E46").Value = '  +3.88%  '

# Row 47
$ws.Range("E47").Value = '  +1.32%  '

# Row 48
$ws.Range("D48").Value = '1.940.24'
$ws.Range("E48").Value = '  -1.40%  '

# Row 49
$ws.Range("E49").Value = '  +0.31%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '112.85'
$ws.Range("E50").Value = '  -0.97%  '

# Row 51
$ws.Range("E51").Value = '  +1.79%  '
